$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Rename E1 "NA_acresk" -> "NAland_acresk", and insert a new column J
# ("lccNA_acresk") before the old lcc columns, shifting old J:U to K:V.
$headers = @("year", "CRPland_acresk", "Cropland_acresk", "Forestland_acresk", "NAland_acresk", "Otherland_acresk", "Pastureland_acresk", "Rangeland_acresk", "Urbanland_acresk", "lccNA_acresk", "lccL1_acresk", "lccL2_acresk", "lccL3_acresk", "lccL4_acresk", "lccL5_acresk", "lccL6_acresk", "lccL7_acresk", "lccL8_acresk", "lccL12_acresk", "lccL34_acresk", "lccL56_acresk", "lccL78_acresk")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data rows (2-8) ---
# Columns A-I keep their meaning (year + landuse areas), with revised values.
# New column J holds the lcc "NA" total; old columns J:U (lcc acreages) shift to K:V.
$years = @(1982, 1987, 1992, 1997, 2002, 2007, 2012)
$data = @(
    @(0.0, 366293.8001255244, 408769.8004858643, 397162.801202178, 112054.60071638972, 183682.80005792528, 417288.50006688386, 50395.500235527754, 518365.7019867152, 29415.90003991127, 285804.700275898, 287258.20012904704, 201368.2001029849, 34436.700000435114, 280848.30010822415, 274476.900233753, 23673.20001332462, 315220.60031580925, 488626.40023203194, 315285.00010865927, 298150.10024707764),
    @(13776.79998434335, 347130.9001113251, 410698.3004591614, 397583.6012101248, 113429.90072029084, 184756.10005189478, 411972.8000700325, 56299.40028312057, 525755.8020436168, 29203.800039298832, 283727.8002592698, 285858.70011573285, 200138.30009755492, 34214.69999996573, 279714.0000971332, 273535.6002253443, 23499.100012376904, 312931.6002985686, 485997.00021328777, 313928.70009709895, 297034.7002377212),
    @(34028.89998526126, 326180.8000878319, 410724.0004300028, 399704.90122456104, 113724.400733307, 179640.80003722757, 408321.6000501439, 63322.40034195781, 534955.8021263555, 28958.700038038194, 281252.1002401337, 283776.60009515285, 198597.60008523613, 33955.29999845475, 278276.5000895187, 272506.5002051294, 23368.700012274086, 310210.8002781719, 482374.200180389, 312231.8000879735, 295875.2002174035),
    @(32694.799986936152, 318596.9000764787, 411963.7003931999, 400059.90122722834, 114703.30073112994, 176979.79999534786, 406921.40004363656, 73728.0004363358, 546260.8022319004, 28614.100037030876, 278231.90020880103, 281246.70006889105, 196635.20006649196, 33697.599996343255, 276691.1000781059, 271132.90019249916, 23137.500010229647, 306846.0002458319, 477881.900135383, 310388.7000744492, 294270.4002027288),
    @(31479.29997756332, 303946.700058423, 412413.7003421709, 401609.20124524087, 115343.10073465854, 182246.49998190254, 406378.500021331, 82230.80052900314, 557037.0023531839, 28167.500035747886, 275477.2001818344, 278901.5000426769, 195199.50004532933, 33458.69999309629, 275395.9000670463, 268967.10016188025, 23043.400009498, 303644.7002175823, 474101.00008800626, 308854.6000601426, 292010.50017137825),
    @(32578.79997328669, 298823.5000416711, 412162.4003028646, 402130.5012490079, 117217.0007352382, 179396.2999684438, 405568.40001321584, 87770.90060656518, 563588.9024424031, 27857.10003284365, 273569.40015506, 277538.1000245139, 194130.30003011227, 33300.2999914065, 274717.9000524655, 267977.7001529038, 22968.10000858456, 301426.50018790364, 471668.40005462617, 308018.200043872, 290945.80016148835),
    @(23949.599979385734, 307765.800039202, 412705.6002696529, 402616.90125477314, 118212.60073700547, 175692.29995437711, 404044.60000356287, 90660.4006523341, 567276.1024981067, 27692.400029584765, 272671.80014347285, 276709.9000074491, 193570.90002006292, 33232.79999014735, 274218.80004697293, 267305.1001460627, 22970.000008434057, 300364.2001730576, 470280.800027512, 307451.6000371203, 290275.10015449673),
)

for ($r = 0; $r -lt $years.Length; $r++) {
    $excelRow = $r + 2
    $ws.Cells.Item($excelRow, 1).Value = $years[$r]
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 2).Value = $rowVals[$c]
    }
}

# New column V falls outside the original A1:U8 used range, so it does not
# inherit the existing numeric style automatically. Match the "0"-format style
# used by the rest of the data cells (same as copying format from column U).
$ws.Range("V2:V8").NumberFormat = $ws.Range("U2").NumberFormat

Write-Host "Updated landu_lcc_totalarea sheet: renamed NA column, added lccNA column, refreshed values"